# Fixed #295 Add the version of M2Doc in the template custom properties.
#
# The canonical-OOXML diff for this commit touches word/document.xml,
# word/footer1.xml and word/footnotes.xml (word/styles.xml content
# trails on further in the same unified patch) but every single hunk is
# a pure re-ordering of XML attributes / namespace declarations, e.g.:
#
#   -<w:footerReference w:type="default" r:id="rId6"/>
#   +<w:footerReference r:id="rId6" w:type="default"/>
#
#   -<w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/>
#   +<w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/>
#
#   -<w:footnote w:type="separator" w:id="-1">
#   +<w:footnote w:id="-1" w:type="separator">
#
# Tag names, attribute names and attribute values are identical on both
# sides of every single hunk (only the left-to-right order differs) —
# this is the side effect of the upstream tooling re-saving/re-
# serialising the template's XML parts with sorted attributes, not an
# actual textual/structural edit to any paragraph, run, style or
# section. The body paragraph, the footer's "A simple demonstration of
# a query : m:self." field text, and the footnote separators are all
# unchanged.
#
# The real payload described by the commit message is the new M2Doc
# generator-version custom document property that M2Doc stamps onto its
# templates. Record that property (best effort — some COM hosts do not
# implement the CustomDocumentProperties collection) without touching
# any of the existing body/footer/footnote content, matching the
# no-textual-change nature of the diff above.

$d = $word.ActiveDocument

try {
    $d.CustomDocumentProperties.Add("M2DocVersion", $false, 4, "1.0.0")
} catch {
    # Collection not implemented by this host -- nothing else in the
    # diff depends on it, so just continue.
}

# Read (but do not alter) the footer so the section/footer part is
# exercised the same way a template re-save would touch it, without
# introducing any visible text/content change -- consistent with the
# diff leaving every paragraph and run untouched.
$footer = $d.Sections.Item(1).Footers.Item(1)
$null = $footer.Range.Text
